$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 764, pushing old rows 764-805 down to 765-806
# (dimension grows from A1:D805 to A1:D806).
$ws.Rows("764:764").Insert()

# Fill in the new row's data: 2026/02/03, 火, 16, 201
#
# Column A holds a date-like string ("2026/02/03"). Assigning it straight to
# .Value would make Excel's COM layer auto-convert it into a date serial
# number (and stamp a date NumberFormat on the cell), but in this workbook
# dates are kept as plain text. Route the text through a formula and then
# flatten the formula to a static value with Copy/PasteSpecial (values-only)
# so the cell ends up as a plain text cell with no extra style applied -
# matching every other date cell in the column.
$ws.Range("A764").Formula = '="2026/02/03"'
$ws.Range("A764").Copy()
$ws.Range("A764").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B764").Value = "火"
$ws.Range("C764").Value = 16
$ws.Range("D764").Value = 201
